$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest forecast row (old row 2); this shifts every subsequent
# row up by one so dates/years in columns A/B/D realign automatically.
$ws.Rows.Item(2).Delete()

# Clear stale forecast residuals (columns C/E) that shifted into rows which
# should now be blank in those columns.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()

# Write the recomputed forecast residual values (columns C/E) for each row.
$ws.Range("C5").Value = -1.120856461469888
$ws.Range("C7").Value = 0.2946875655135184
$ws.Range("C9").Value = 1.242549344471122
$ws.Range("C11").Value = 1.180122656701199
$ws.Range("E11").Value = 0.9288717675470126
$ws.Range("C12").Value = 0.9006569003772169
$ws.Range("E12").Value = 0.4141666450523163
$ws.Range("C13").Value = 1.381744454873757
$ws.Range("E13").Value = 1.416624765035412
$ws.Range("C14").Value = 1.298949644080372
$ws.Range("E14").Value = 0.7446484179501223
$ws.Range("C15").Value = 1.404039427736437
$ws.Range("E15").Value = 1.075154359849861
$ws.Range("C16").Value = 1.485127130420993
$ws.Range("E16").Value = 0.8988967199517361
$ws.Range("C17").Value = 1.593617458167307
$ws.Range("E17").Value = 1.214249019249602
$ws.Range("C18").Value = 1.938044824544427
$ws.Range("E18").Value = 1.197301207077017
$ws.Range("C19").Value = 1.885212754467758
$ws.Range("E19").Value = 1.180565832117297
$ws.Range("C20").Value = 1.730502563828185
$ws.Range("E20").Value = 1.20920901052266
$ws.Range("C21").Value = 1.877372574190161
$ws.Range("E21").Value = 1.638669199130427
$ws.Range("C22").Value = 1.95207711838874
$ws.Range("E22").Value = 1.265658248521984
$ws.Range("C23").Value = 2.211325510218898
$ws.Range("E23").Value = 1.513838358900466
$ws.Range("C24").Value = 2.320468093590722
$ws.Range("E24").Value = 1.753251432658076
$ws.Range("C25").Value = 2.279995067217899
$ws.Range("E25").Value = 1.540918326052476
$ws.Range("C26").Value = 1.036239343320755
$ws.Range("E26").Value = 1.010297891741785
$ws.Range("C27").Value = 1.095903126316466
$ws.Range("E27").Value = 1.063472944477306
$ws.Range("C28").Value = 1.142800289129831
$ws.Range("E28").Value = 1.156859643506358
$ws.Range("C29").Value = 1.025257057800411
$ws.Range("E29").Value = 0.5773070399857971
$ws.Range("C30").Value = 0.2803494251802263
$ws.Range("E30").Value = 0.855281508838468
$ws.Range("C31").Value = 0.8137456736830195
$ws.Range("E31").Value = 1.30966355756772
$ws.Range("C32").Value = -1.788000783651811
$ws.Range("E32").Value = -4.916332952555802
$ws.Range("C33").Value = -1.788000783651811
$ws.Range("E33").Value = -0.415982961498651
$ws.Range("C34").Value = -1.484481523646708
$ws.Range("E34").Value = 0.7587717871427202
$ws.Range("C35").Value = -1.434438137829841
$ws.Range("E35").Value = 0.8159375071586261
$ws.Range("C36").Value = -1.098964423305859
$ws.Range("E36").Value = 1.122873944479474
$ws.Range("C37").Value = -1.098964423305859
$ws.Range("E37").Value = 1.426719405738508
$ws.Range("C38").Value = 1.659091542859148
$ws.Range("E38").Value = 1.163281843182573
$ws.Range("C39").Value = 1.85385197842538
$ws.Range("E39").Value = 1.2808239555127
$ws.Range("C40").Value = 1.916393754370604
$ws.Range("E40").Value = 1.358354219817404
$ws.Range("C41").Value = 1.916393754370604
$ws.Range("E41").Value = 0.5205511175203181
$ws.Range("C42").Value = -0.2882829247660479
$ws.Range("E42").Value = 0.9814522951840488
$ws.Range("C43").Value = -0.6079479926716203
$ws.Range("E43").Value = 0.8021760422591839
$ws.Range("C44").Value = -0.7359525160776204
$ws.Range("E44").Value = 0.6489285084265051
$ws.Range("C45").Value = -0.7359525160776204
$ws.Range("E45").Value = 0.3517304536567734
$ws.Range("C46").Value = 0.1328794705491632
$ws.Range("E46").Value = 0.8679080298191755
$ws.Range("C47").Value = -0.06520462171909491
$ws.Range("E47").Value = 0.7367476213790747
$ws.Range("C48").Value = -0.187152549496028
$ws.Range("E48").Value = 0.5818712200032161
$ws.Range("C49").Value = -0.187152549496028
$ws.Range("E49").Value = 0.4186921370205043
$ws.Range("C50").Value = 0.5283432505880592
$ws.Range("E50").Value = 0.8122218530175696
$ws.Range("C51").Value = 0.5869668956646645
$ws.Range("E51").Value = 0.8208952814083625
$ws.Range("C52").Value = 0.6150340712028246
$ws.Range("E52").Value = 0.831632580682462

